$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 23-24 used to be a merged continuation (B23:B24) duplicating the
# "Vertical dilution of precision" block. Replace with new RTK/GPS-advanced-fix
# info and drop the merge.
$ws.Range("B23:B24").UnMerge()

$ws.Range("A23").Value = $null
$ws.Range("B23").Value = "GPS advanced fix"
$ws.Range("C23").Value = "N/A"
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = "0: no advanced fix, 1: DGPS, 2: RTK_FLOAT, 3: RTK_FIXED"

$ws.Range("A24").Value = $null
$ws.Range("B24").Value = "RESERVED"
$ws.Range("C24").Value = $null
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = $null
